# Add the missing "explanation about solution" textbox to the
# "Steps.ScalingOut" slide (the slide with Title 1 / Picture 4 / TextBox 2 /
# Table 5), right after the existing table, matching the author's commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Shape.Left/Top/Width/Height (and AddTextbox's coordinates) are expressed
# in points in the PowerPoint object model, while the target OOXML values
# below are EMUs (1 pt = 12700 EMU) -- convert down before creating/sizing.
$emuPerPt = 12700
$left   = 654340   / $emuPerPt
$top    = 5046983  / $emuPerPt
$width  = 10704353 / $emuPerPt
$height = 923330   / $emuPerPt

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 3"

# spPr: no fill, word-wrapped + shrink-to-fit text box (matches the other
# free-floating textboxes already on this slide deck).
$shp.Fill.Visible = $false
$shp.TextFrame.WordWrap = -1
$shp.TextFrame.AutoSize = 1

$tr = $shp.TextFrame.TextRange
$tr.Text = "The issue when several "
$tr = $tr.InsertAfter("SignUp")
$tr = $tr.InsertAfter(" Functions are booking the last place is solved by predefining Places in ")
$tr = $tr.InsertAfter("CoursePlaces")
$tr = $tr.InsertAfter(" table. Then instead of insert operation we use update operation. Update operation uses optimistic locking with ")
$tr = $tr.InsertAfter("Etags")
$tr = $tr.InsertAfter(". So we delegate responsibility of solving this concurrency issue to Azure Table Storage.")
